$d = $word.ActiveDocument

# Find and replace "50" with "25" in the points paragraph
$d.Content.Find.Execute("50", $true, $false, $false, $false, $false,
                         $true, 1, $false, "25", 2)
